# Jogos_da_Semana_FlashScore_2024-11-09.xlsx refresh (odds/fixtures update).
#
# * The old row 13 fixture (USA - USL CHAMPIONSHIP, Louisville City vs Rhode
#   Island) is dropped, so the sheet shrinks from A1:BD13 to A1:BD12.
# * Row 4 is replaced with a new fixture (COLOMBIA - PRIMERA A, Bucaramanga
#   vs Fortaleza).
# * Rows 11/12 now hold the fixtures that used to be in rows 12/13 (Inter
#   Miami vs Atlanta Utd, and Louisville City vs Rhode Island), with refreshed
#   odds.
# * Rows 2, 3 and 10 keep their fixtures but get updated odds in a handful of
#   columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last fixture row - everything above keeps its row number.
$ws.Rows("13").Delete()

# --- Row 2 ---
$ws.Range("G2").Value = 2.9
$ws.Range("I2").Value = 3.1
$ws.Range("L2").Value = 4.33
$ws.Range("M2").Value = 1.22
$ws.Range("N2").Value = 4
$ws.Range("Q2").Value = 4.2
$ws.Range("R2").Value = 1.22
$ws.Range("S2").Value = 1.95
$ws.Range("T2").Value = 1.85
$ws.Range("W2").Value = 5
$ws.Range("Z2").Value = 34
$ws.Range("AC2").Value = 4
$ws.Range("AI2").Value = 12
$ws.Range("AJ2").Value = 15
$ws.Range("AO2").Value = 23
$ws.Range("AT2").Value = 1.83
$ws.Range("AW2").Value = 4.5
$ws.Range("AX2").Value = 23

# --- Row 3 ---
$ws.Range("N3").Value = 9
$ws.Range("Y3").Value = 9
$ws.Range("AJ3").Value = 19
$ws.Range("BA3").Value = 151

# --- Row 4 ---
$ws.Range("A4").Value = "nwbcDY6U"
$ws.Range("C4").Value = "22:20"
$ws.Range("D4").Value = "COLOMBIA - PRIMERA A"
$ws.Range("E4").Value = "Bucaramanga"
$ws.Range("F4").Value = "Fortaleza"
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 3.9
$ws.Range("J4").Value = 2.75
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 4.75
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.38
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 8.5
$ws.Range("Y4").Value = 9.5
$ws.Range("Z4").Value = 17
$ws.Range("AA4").Value = 19
$ws.Range("AB4").Value = 34
$ws.Range("AC4").Value = 7
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 19
$ws.Range("AF4").Value = 67
$ws.Range("AG4").Value = 201
$ws.Range("AH4").Value = 9
$ws.Range("AI4").Value = 19
$ws.Range("AJ4").Value = 15
$ws.Range("AL4").Value = 41
$ws.Range("AM4").Value = 41
$ws.Range("AN4").Value = 3.75
$ws.Range("AO4").Value = 12
$ws.Range("AP4").Value = 26
$ws.Range("AQ4").Value = 41
$ws.Range("AR4").Value = 67
$ws.Range("AS4").Value = 251
$ws.Range("AT4").Value = 2.38
$ws.Range("AU4").Value = 9
$ws.Range("AV4").Value = 67
$ws.Range("AW4").Value = 5.5
$ws.Range("AY4").Value = 34
$ws.Range("AZ4").Value = 81
$ws.Range("BA4").Value = 126
$ws.Range("BB4").Value = 351
$ws.Range("BC4").Value = 126
$ws.Range("BD4").Value = 126

# --- Row 10 ---
$ws.Range("O10").Value = 1.36
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = 2.15
$ws.Range("R10").Value = 1.67
$ws.Range("S10").Value = 1.41
$ws.Range("T10").Value = 2.62

# --- Row 11 ---
$ws.Range("A11").Value = "634HgM7l"
$ws.Range("C11").Value = "22:00"
$ws.Range("E11").Value = "Inter Miami"
$ws.Range("F11").Value = "Atlanta Utd"
$ws.Range("G11").Value = 1.36
$ws.Range("H11").Value = 5.25
$ws.Range("I11").Value = 7.5
$ws.Range("J11").Value = 1.73
$ws.Range("K11").Value = 2.88
$ws.Range("L11").Value = 6
$ws.Range("M11").Value = 1.01
$ws.Range("N11").Value = 23
$ws.Range("O11").Value = 1.08
$ws.Range("P11").Value = 8
$ws.Range("Q11").Value = 1.33
$ws.Range("R11").Value = 3.4
$ws.Range("S11").Value = 1.17
$ws.Range("T11").Value = 4.33
$ws.Range("U11").Value = 1.57
$ws.Range("V11").Value = 2.25
$ws.Range("W11").Value = 12
$ws.Range("X11").Value = 9.5
$ws.Range("Y11").Value = 9
$ws.Range("Z11").Value = 11
$ws.Range("AA11").Value = 10
$ws.Range("AB11").Value = 19
$ws.Range("AC11").Value = 23
$ws.Range("AD11").Value = 11
$ws.Range("AF11").Value = 41
$ws.Range("AG11").Value = 126
$ws.Range("AH11").Value = 26
$ws.Range("AI11").Value = 41
$ws.Range("AJ11").Value = 21
$ws.Range("AK11").Value = 81
$ws.Range("AL11").Value = 41
$ws.Range("AO11").Value = 6.5
$ws.Range("AP11").Value = 13
$ws.Range("AQ11").Value = 15
$ws.Range("AR11").Value = 29
$ws.Range("AS11").Value = 67
$ws.Range("AT11").Value = 4.33
$ws.Range("AV11").Value = 41
$ws.Range("AW11").Value = 9
$ws.Range("AX11").Value = 34
$ws.Range("AZ11").Value = 101
$ws.Range("BB11").Value = 151
$ws.Range("BC11").Value = 301
$ws.Range("BD11").Value = 176

# --- Row 12 ---
$ws.Range("A12").Value = "UJd8iehn"
$ws.Range("C12").Value = "21:30"
$ws.Range("D12").Value = "USA - USL CHAMPIONSHIP"
$ws.Range("E12").Value = "Louisville City"
$ws.Range("F12").Value = "Rhode Island"
$ws.Range("G12").Value = 1.6
$ws.Range("H12").Value = 4.1
$ws.Range("I12").Value = 4.4
$ws.Range("J12").Value = 2.07
$ws.Range("K12").Value = 2.47
$ws.Range("L12").Value = 4.55
$ws.Range("M12").Value = 1.03
$ws.Range("N12").Value = 9.25
$ws.Range("O12").Value = 1.16
$ws.Range("P12").Value = 4.5
$ws.Range("Q12").Value = 1.5
$ws.Range("R12").Value = 2.4
$ws.Range("S12").Value = 1.27
$ws.Range("T12").Value = 3.4
$ws.Range("W12").Value = 9.75
$ws.Range("X12").Value = 9.25
$ws.Range("Y12").Value = 8.25
$ws.Range("Z12").Value = 13
$ws.Range("AA12").Value = 11.5
$ws.Range("AC12").Value = 9.25
$ws.Range("AD12").Value = 8.5
$ws.Range("AE12").Value = 14
$ws.Range("AF12").Value = 45
$ws.Range("AG12").Value = 250
$ws.Range("AH12").Value = 17.5
$ws.Range("AI12").Value = 30
$ws.Range("AJ12").Value = 15
$ws.Range("AK12").Value = 75
$ws.Range("AL12").Value = 37
$ws.Range("AM12").Value = 35
$ws.Range("AN12").Value = 3.8
$ws.Range("AO12").Value = 7.5
$ws.Range("AP12").Value = 14
$ws.Range("AQ12").Value = 21
$ws.Range("AR12").Value = 40
$ws.Range("AS12").Value = 150
$ws.Range("AT12").Value = 3.4
$ws.Range("AU12").Value = 7
$ws.Range("AV12").Value = 50
$ws.Range("AW12").Value = 6.7
$ws.Range("AX12").Value = 23
$ws.Range("AY12").Value = 25
$ws.Range("AZ12").Value = 120
$ws.Range("BA12").Value = 120
$ws.Range("BB12").Value = 250
$ws.Range("BC12").Value = 51
$ws.Range("BD12").Value = 51
